# Edit: insert two new price-record rows (Cebollín, Femacal de La Calera)
# right before the existing row 150, shifting the remaining historical
# rows down by two. The two new rows hold a fresh "Primera"/"Segunda"
# quote pair for date serial 44438 (2021-08-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 150; existing rows 150-254 shift to 152-256
$ws.Rows.Item(150).Insert()
$ws.Rows.Item(150).Insert()

# --- New row 150 (Primera) ---
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(150, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44438
$ws.Cells.Item(150, 5).Value = 5
$ws.Cells.Item(150, 6).Value = 100112037
$ws.Cells.Item(150, 7).Value = "Cebollín"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 230
$ws.Cells.Item(150, 11).Value = 4000
$ws.Cells.Item(150, 12).Value = 4000
$ws.Cells.Item(150, 13).Value = 4000
$ws.Cells.Item(150, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(150, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(150, 16).Value = 111
$ws.Cells.Item(150, 17).Value = 36
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# --- New row 151 (Segunda) ---
$ws.Cells.Item(151, 1).Value = 3
$ws.Cells.Item(151, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44438
$ws.Cells.Item(151, 5).Value = 5
$ws.Cells.Item(151, 6).Value = 100112037
$ws.Cells.Item(151, 7).Value = "Cebollín"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Segunda"
$ws.Cells.Item(151, 10).Value = 220
$ws.Cells.Item(151, 11).Value = 2500
$ws.Cells.Item(151, 12).Value = 2500
$ws.Cells.Item(151, 13).Value = 2500
$ws.Cells.Item(151, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(151, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(151, 16).Value = 69
$ws.Cells.Item(151, 17).Value = 36
$ws.Cells.Item(151, 18).Value = "Hortaliza"

# Ensure the date cells keep the same date number format as the rest of column D
$ws.Cells.Item(150, 4).NumberFormat = $ws.Cells.Item(149, 4).NumberFormat
$ws.Cells.Item(151, 4).NumberFormat = $ws.Cells.Item(149, 4).NumberFormat
